# Append new job listing row at the top of the data (row 2) and refresh the
# "取得日時" (fetched-at) timestamp on every existing row, matching a re-scrape
# pass at 2025-09-19 01:43:41. All rows below the header shift down by one and
# keep their original cell order/styling; the URL hyperlinks are rebuilt so the
# rId mapping stays 1:1 with the (now shifted) row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Full target data set for rows 2..20 (A:取得日時, B:タイトル, C:カテゴリ, D:価格,
# E:締切, F:URL, G:優先度スコア, H:スキル概要 -- H left blank ("") where the
# source row has no skill summary, matching the original sheet).
$data = @(
    ('2025-09-19 01:43:41', '【急募】行政情報収集・要約 AI エージェント開発の依頼', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5396253', '368', '🔥AI,Ai ◆開発'),
    ('2025-09-19 01:43:41', '【報酬4万円〜|相談可能】n8n構築者募集|AIワークフロー構築が得意な方を探しています', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5396220', '303', '🔥AI,Ai'),
    ('2025-09-19 01:43:41', 'AIチャットボットのβ版テスト参加者募集!', 'システム開発', '10,000 円 ~ 20,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395790', '295', '🔥AI,Ai'),
    ('2025-09-19 01:43:41', '【急募】JUSTDBとOPERAcloudのAPI連携開発者募集', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5396169', '250', '🔥API ◆開発'),
    ('2025-09-19 01:43:41', '1688アリババの商品情報の抽出のスクレイピングの開発 exe形式の自動ツール', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5387065', '173', '◆ツール,開発'),
    ('2025-09-19 01:43:41', '【急募】入力ミス防止のためのkintone Javascript開発者募集', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395671', '128', '★Java ◆開発'),
    ('2025-09-19 01:43:41', '初回 Laravel Livewireを使ったWebシステム開発の募集', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395502', '125', '◆開発,システム開発'),
    ('2025-09-19 01:43:41', '2026年度新入社員研修Javaサブ講師', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395710', '85', '★Java'),
    ('2025-09-19 01:43:41', '2026年度新入社員研修Javaサブ講師 (4~6月)', 'システム開発', '1,000,000 円 ~ 3,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395713', '85', '★Java'),
    ('2025-09-19 01:43:41', '【急募】iOSアプリのAdMobメディエーション入札接続とeCPM改善', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395931', '38', '◇アプリ'),
    ('2025-09-19 01:43:41', '【急募】HP保守管理とSEO対策の専門家を探しています', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5396003', '33', '◇管理'),
    ('2025-09-19 01:43:41', '【急募】ストアーズ連携LINE予約サイト制作のプロを探しています!', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395825', '33', '◇サイト'),
    ('2025-09-19 01:43:41', '【継続案件|お気軽にご応募ください】WebシステムのQAエンジニア募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395741', '40', ''),
    ('2025-09-19 01:43:41', '社内のFAQの構築', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5396173', '25', ''),
    ('2025-09-19 01:43:41', '〖リモート可〗Delphiエンジニア募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5341051', '25', ''),
    ('2025-09-19 01:43:41', '【急募】ウェブプラットフォームのMVP制作を依頼します!', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5396017', '18', ''),
    ('2025-09-19 01:43:41', 'MT5用EA(ex5)ファイルのデコンパイル', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5396009', '13', ''),
    ('2025-09-19 01:43:41', '【急募】JotformとGoogleスプレッドシート連携のエキスパート募集!', 'システム開発', '5,000 円 ~ 10,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395809', '10', ''),
    ('2025-09-19 01:43:41', '【MT4】ゴールドの取引を行うEAのサンプルソース納品', 'システム開発', '10,000 円 ~ 20,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5395799', '10', '')
)

# Clear every existing hyperlink on the sheet -- deleting any single cell's
# Hyperlinks collection clears the whole worksheet collection in this engine,
# so do it once up front and rebuild all F-column links below in row order.
$ws.Range("F2").Hyperlinks.Delete()

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $ws.Cells.Item($row, 7).Value = [int]$item[6]
    if ($item[7] -ne "") {
        $ws.Cells.Item($row, 8).Value = $item[7]
    }

    $linkCell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($linkCell, $item[5])
    $linkCell.Style = "Hyperlink"

    $row = $row + 1
}

